# BombSquad_UserManual.docx edit script
#
# This script performs three textual edits in the document body:
#
# 1) In the opening "Bomb Squad is a cooperative game..." paragraph,
#    consolidate the several runs that make up the first sentence
#    ("...puzzles on the DE-2 FPGA board with the objective, to defuse
#    "the bomb". ...a bomb ") into a single run. The visible text is
#    unchanged; only the run/character-run boundaries are merged.
#
# 2) In the same paragraph, consolidate the runs describing the two
#    player roles ("...the player performing the actions and is not
#    allowed...manual.") into a single run. Again, the visible text is
#    unchanged; only run boundaries are merged.
#
# 3) In the "Login Sequence for the Bomb Squad Game" section, update the
#    sentence describing the login sequence:
#       "...was equipped with a login verifier sequence designed to
#        identify the user by the input of credentials."
#    becomes
#       "...was equipped with a credential verifier login sequence
#        designed to identify a user by the input of credentials."
#    The `_GoBack` bookmark that sits in the middle of this sentence is
#    left in place (not deleted) by editing the text strictly on either
#    side of it.

$d = $word.ActiveDocument

# --- Locate the two paragraphs we need without disturbing a later,
# --- near-duplicate paragraph (an appendix sample) that must NOT be
# --- touched.
$paras = $d.Paragraphs
$introPara = $null
$loginPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if (($introPara -eq $null) -and ($t -like "*Bomb Squad is a cooperative game*DE-2*")) {
        $introPara = $i
    }
    if (($loginPara -eq $null) -and ($t -like "*Bomb Squa*interactive game w*")) {
        $loginPara = $i
    }
}

# ---------------------------------------------------------------------
# Edit 1: merge the runs of the first sentence of the intro paragraph.
# ---------------------------------------------------------------------
$r1 = $paras.Item($introPara).Range
$find1 = "Bomb Squad is a cooperative game in which the objective is to perform various tasks and solve puzzles on the DE-2 FPGA board with the objective, to defuse " + [char]0x201C + "the bomb" + [char]0x201D + ". The instructions for defusing the bomb will be provided in a bomb "
$null = $r1.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $find1, 2)

# ---------------------------------------------------------------------
# Edit 2: merge the runs of the player-roles sentence of the intro
# paragraph (the run(s) containing "defusal" -- wrapped in proofErr
# spell-check markers -- are left alone, exactly as in the target).
# ---------------------------------------------------------------------
$r2 = $paras.Item($introPara).Range
$find2 = " manual. There are two player roles in this game: the technician and the specialist. The technician is the player performing the actions and is not allowed to look at the manual at any time. The specialist is allowed to look at manual and must verbally communicate with the technician to perform the steps outlined in the manual."
$null = $r2.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $find2, 2)

# ---------------------------------------------------------------------
# Edit 3: update the login-sequence sentence, without disturbing the
# _GoBack bookmark that sits between "...game w" and "as equipped...".
# Handle the text after the bookmark first, then the text before it, so
# neither Find.Execute match ever has to span across the bookmark.
# ---------------------------------------------------------------------

# 3a) After the bookmark: "as equipped with a login verifier sequence
#     designed to identify the user by the input of credentials."
#     -> " user by the input of credentials."
$r3a = $paras.Item($loginPara).Range
$find3a = "as equipped with a login verifier sequence designed to identify the user by the input of credentials."
$repl3a = " user by the input of credentials."
$null = $r3a.Find.Execute($find3a, $true, $false, $false, $false, $false, $true, 1, $false, $repl3a, 2)

# 3b) Before the bookmark: "d interactive game w"
#     -> "d interactive game was equipped with a credential verifier
#         login sequence designed to identify a"
$r3b = $paras.Item($loginPara).Range
$find3b = "d interactive game w"
$repl3b = "d interactive game was equipped with a credential verifier login sequence designed to identify a"
$null = $r3b.Find.Execute($find3b, $true, $false, $false, $false, $false, $true, 1, $false, $repl3b, 2)
